$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (prices and volume % changes).
# Cells in column D that look like plain numbers are forced to text format
# ('@') before assignment so Excel preserves their original textual
# representation (e.g. trailing zeros, multi-dot formatted prices) instead of
# silently converting them to numeric values.

$ws.Range('D2').Value = '68.888.54'
$ws.Range('E2').Value = '  -3.02%  '
$ws.Range('D3').Value = '3.438.95'
$ws.Range('E3').Value = '  -5.75%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '573.93'
$ws.Range('E5').Value = '  -5.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '188.90'
$ws.Range('E6').Value = '  -5.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.607'
$ws.Range('E7').Value = '  -3.51%  '
$ws.Range('D8').Value = '3.428.97'
$ws.Range('E8').Value = '  -5.47%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.203'
$ws.Range('E10').Value = '  -7.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.612'
$ws.Range('E11').Value = '  -5.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '50.88'
$ws.Range('E12').Value = '  -5.74%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000282'
$ws.Range('E13').Value = '  -8.25%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.02'
$ws.Range('E14').Value = '  -5.85%  '
$ws.Range('D15').Value = '3.982.51'
$ws.Range('E15').Value = '  -5.49%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '635.35'
$ws.Range('E16').Value = '  +4.71%  '
$ws.Range('D17').Value = '68.693.14'
$ws.Range('E17').Value = '  -3.31%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.432.61'
$ws.Range('E18').Value = '  -5.43%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.120'
$ws.Range('E19').Value = '  -2.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.13'
$ws.Range('E20').Value = '  -7.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.98'
$ws.Range('E21').Value = '  -5.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.934'
$ws.Range('E22').Value = '  -6.72%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.62'
$ws.Range('E23').Value = '  -3.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.30'
$ws.Range('E24').Value = '  +0.09%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '97.96'
$ws.Range('E25').Value = '  -5.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.24'
$ws.Range('E26').Value = '  -8.57%  '
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.82'
$ws.Range('E27').Value = '  -6.29%  '
$ws.Range('B28').Value = 'LEO'
$ws.Range('C28').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.07'
$ws.Range('E28').Value = '  +1.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.79'
$ws.Range('E29').Value = '  -7.99%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.14'
$ws.Range('E30').Value = '  -6.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '32.06'
$ws.Range('E31').Value = '  -5.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.29'
$ws.Range('E32').Value = '  -9.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.63'
$ws.Range('E33').Value = '  -8.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '11.46'
$ws.Range('E34').Value = '  -7.01%  '
$ws.Range('E35').Value = '  -4.31%  '
$ws.Range('E36').Value = '  -7.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  -0.21%  '
$ws.Range('D38').Value = '3.639.17'
$ws.Range('E38').Value = '  -8.82%  '
$ws.Range('D39').Value = '0.0₃0782'
$ws.Range('E39').Value = '  -12.47%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '489.16'
$ws.Range('E40').Value = '  -5.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.87'
$ws.Range('E41').Value = '  -7.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.42'
$ws.Range('E42').Value = '  -4.00%  '
$ws.Range('B43').Value = 'CoreDAO'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.52'
$ws.Range('E43').Value = '  +70.51%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.365'
$ws.Range('E44').Value = '  -6.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.132'
$ws.Range('E45').Value = '  -3.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '33.91'
$ws.Range('E46').Value = '  -7.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0434'
$ws.Range('E47').Value = '  -6.40%  '
$ws.Range('E48').Value = '  -4.19%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.78'
$ws.Range('E49').Value = '  -5.41%  '
$ws.Range('E50').Value = '  -5.54%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.998'
$ws.Range('E51').Value = '  -0.33%  '
